$d = $word.ActiveDocument

# The two checkbox content controls ("Branch created for analysis..." and
# "Use GitHub...") need to flip from unchecked to checked.
foreach ($cc in $d.ContentControls) {
    if ($cc.ID -eq -916699948 -or $cc.ID -eq 2108696346) {
        $cc.Checked = $true
        $cc.Range.Text = [char]0x2612
    }
}
